# Append the 09/10/2025 profit figure as a new row (row 24) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(24, 1)
$profitCell = $ws.Cells.Item(24, 2)

# Force the date column to be stored as literal text (matching the existing
# "MM/DD/YYYY" string entries above it) instead of Excel's automatic date
# recognition/serial-number conversion.
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/10/2025"
$profitCell.Value = 14898.33

# Re-apply the plain formatting used by the other date cells in column A so
# the new cell doesn't end up with a one-off "text" number format applied to
# it (keeps it visually/structurally consistent with A2:A23).
$ws.Range("A2").Copy()
$ws.Range("A24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
